# Trade #8 closed at 2026-02-17 23:52:48 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1500.13   # Current Capital
$wsSummary.Range("B4").Value = 0.13      # Total P&L $
$wsSummary.Range("B5").Value = 0.33      # Total P&L %
$wsSummary.Range("B6").Value = 8         # Total Trades
$wsSummary.Range("B7").Value = 4         # Winning Trades
$wsSummary.Range("B9").Value = 50        # Win Rate %

# --- Strategy Status sheet (MarketMaking row) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C6").Value = 100.13
$wsStatus.Range("D6").Value = 8
$wsStatus.Range("E6").Value = 0.13
$wsStatus.Range("F6").Value = 0.13
$wsStatus.Range("G6").Value = 50

# --- All Trades & MarketMaking sheets (Trade #8 row, row 9) ---
$tradeSheets = @("All Trades", "MarketMaking")
foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("G9").Value = 0.25
    $ws.Range("H9").Value = "CLOSED"
    $ws.Range("I9").Value = 39.1176
    $ws.Range("J9").Value = 0.07000000000000001
    $ws.Range("K9").Value = 100.13
    $ws.Range("P9").Value = "early_exit"
    $ws.Range("Q9").Value = 0.13
}
